# Lab02.pptx / week02 edit:
#   "    char a = 127;"  ->  "    signed char a = 127;"
# The author split the single run into two runs:
#   run 1 (new):      "    signed char "
#   run 2 (original): "a = 127;"   (keeps its original rPr incl. dirty="0")
#
# The code listing lives in slide 11, in the third shape ("文本框 6" / "Text Box 6").

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(11)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

$fullText = $tr.Text
$needle   = "    char a = 127;"
$pos      = $fullText.IndexOf($needle)

if ($pos -lt 0) {
    throw "Could not locate target line '    char a = 127;' in shape text"
}

# PowerPoint TextRange indices are 1-based.
$start = $pos + 1
$len   = 9   # "    char " (four leading spaces + "char ") -> becomes "    signed char "

$target = $tr.Characters($start, $len)
$target.Text = "    signed char "
